$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3480920195579529
$ws.Range("B1").Value = 2.251837968826294
$ws.Range("C1").Value = 4.798030853271484
$ws.Range("D1").Value = 1.721617937088013
$ws.Range("E1").Value = 0.8628863096237183
